$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(12, 4).Value = 98.04000000000001
$ws.Cells.Item(12, 5).Value = 4.201
$ws.Cells.Item(12, 6).Value = 3.61
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(24, 4).Value = 98.04000000000001
$ws.Cells.Item(24, 5).Value = 4.58
$ws.Cells.Item(24, 6).Value = 3.61
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(35, 5).Value = 0
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(36, 4).Value = 98.04000000000001
$ws.Cells.Item(36, 5).Value = 0
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 4.32
$ws.Cells.Item(37, 4).Value = 0
$ws.Cells.Item(37, 5).Value = 0
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 7).Value = 0
$ws.Cells.Item(38, 5).Value = 0
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(39, 5).Value = 0
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(40, 5).Value = 0
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(41, 5).Value = 0
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(43, 5).Value = 0
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(45, 5).Value = 0
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(46, 5).Value = 0
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(47, 6).Value = 0
$ws.Cells.Item(48, 4).Value = 98.04000000000001
$ws.Cells.Item(48, 5).Value = 0
$ws.Cells.Item(48, 6).Value = 0
$ws.Cells.Item(48, 7).Value = 4.05
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(49, 5).Value = 0
$ws.Cells.Item(49, 6).Value = 0
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(52, 5).Value = 0
$ws.Cells.Item(52, 6).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(54, 5).Value = 0
$ws.Cells.Item(54, 6).Value = 0
$ws.Cells.Item(55, 5).Value = 0
$ws.Cells.Item(55, 6).Value = 0
$ws.Cells.Item(56, 5).Value = 0
$ws.Cells.Item(56, 6).Value = 0
$ws.Cells.Item(57, 5).Value = 0
$ws.Cells.Item(57, 6).Value = 0
$ws.Cells.Item(58, 5).Value = 0
$ws.Cells.Item(58, 6).Value = 0
$ws.Cells.Item(59, 5).Value = 0
$ws.Cells.Item(59, 6).Value = 0
$ws.Cells.Item(60, 4).Value = 98.04000000000001
$ws.Cells.Item(60, 5).Value = 0
$ws.Cells.Item(60, 6).Value = 0
$ws.Cells.Item(60, 7).Value = 3.95
$ws.Cells.Item(61, 4).Value = 0
$ws.Cells.Item(61, 5).Value = 0
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(62, 5).Value = 0
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(66, 5).Value = 0
$ws.Cells.Item(66, 6).Value = 0
$ws.Cells.Item(67, 5).Value = 0
$ws.Cells.Item(67, 6).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(72, 4).Value = 98.04000000000001
$ws.Cells.Item(72, 5).Value = 0
$ws.Cells.Item(72, 6).Value = 0
$ws.Cells.Item(72, 7).Value = 4.01
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 0
